$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 642.82355
$ws.Range("I92").Value = 517.0909
$ws.Range("J92").Value = 873.3333
$ws.Range("K92").Value = 517.0909
$ws.Range("L92").Value = 873.3333
$ws.Range("M92").Value = 730.9091
$ws.Range("N92").Value = -3369.3333

$ws.Range("H95").Value = 86713.28999999999
$ws.Range("J95").Value = 86713.28999999999
$ws.Range("L95").Value = 86713.28999999999
$ws.Range("N95").Value = -92205.28999999999

$ws.Range("H113").Value = 267212.16
$ws.Range("I113").Value = 349797.2
$ws.Range("J113").Value = 2940
$ws.Range("K113").Value = 349797.2
$ws.Range("L113").Value = 2940
$ws.Range("M113").Value = -346543.2
$ws.Range("N113").Value = -9448

$ws.Range("H132").Value = 1374616.2
$ws.Range("I132").Value = 1488934.4
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 4466803.199999999
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -4464273.199999999
$ws.Range("N132").Value = -13460

$ws.Range("H137").Value = 1859.25
$ws.Range("I137").Value = 949.375
$ws.Range("J137").Value = 2769.125
$ws.Range("K137").Value = 2848.125
$ws.Range("L137").Value = 8307.375
$ws.Range("M137").Value = -298.125
$ws.Range("N137").Value = -13407.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7419.547
$ws.Range("I32").Value = 4560.931
$ws.Range("K32").Value = 4560.931
$ws.Range("M32").Value = -4273.931

$ws.Range("H74").Value = 805.94446
$ws.Range("J74").Value = 1027.5454
$ws.Range("L74").Value = 1027.5454
$ws.Range("N74").Value = -2775.5454

$ws.Range("H77").Value = 805.94446
$ws.Range("J77").Value = 1027.5454
$ws.Range("L77").Value = 5137.727
$ws.Range("N77").Value = -13873.727

$ws.Range("H132").Value = 1381.091
$ws.Range("I132").Value = 1229.7778
$ws.Range("K132").Value = 3689.3334
$ws.Range("M132").Value = -1159.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 23260358
$ws.Range("I20").Value = 43484692
$ws.Range("J20").Value = 2372.5
$ws.Range("K20").Value = 43484692
$ws.Range("L20").Value = 2372.5
$ws.Range("M20").Value = -43484445
$ws.Range("N20").Value = -2866.5

$ws.Range("H86").Value = 1662
$ws.Range("I86").Value = 1496.25
$ws.Range("J86").Value = 1851.4286
$ws.Range("K86").Value = 1496.25
$ws.Range("L86").Value = 1851.4286
$ws.Range("M86").Value = -373.25
$ws.Range("N86").Value = -4097.4286

$ws.Range("H89").Value = 1662
$ws.Range("I89").Value = 1496.25
$ws.Range("J89").Value = 1851.4286
$ws.Range("K89").Value = 7481.25
$ws.Range("L89").Value = 9257.143
$ws.Range("M89").Value = -1865.25
$ws.Range("N89").Value = -20489.143

$ws.Range("H105").Value = 1264218.5
$ws.Range("I105").Value = 2274207.2
$ws.Range("J105").Value = 1732.625
$ws.Range("K105").Value = 2274207.2
$ws.Range("L105").Value = 1732.625
$ws.Range("M105").Value = -2272460.2
$ws.Range("N105").Value = -5226.625

$ws.Range("H134").Value = 1404.1132
$ws.Range("I134").Value = 1039.2778
$ws.Range("J134").Value = 2176.7058
$ws.Range("K134").Value = 3117.8334
$ws.Range("L134").Value = 6530.117400000001
$ws.Range("M134").Value = -582.8334000000004
$ws.Range("N134").Value = -11600.1174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2158.074
$ws.Range("I99").Value = 1875.2222
$ws.Range("J99").Value = 2723.7778
$ws.Range("K99").Value = 1875.2222
$ws.Range("L99").Value = 2723.7778
$ws.Range("M99").Value = -377.2221999999999
$ws.Range("N99").Value = -5719.7778

$ws.Range("H105").Value = 6100
$ws.Range("I105").Value = 7120
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 7120
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = -5373
$ws.Range("N105").Value = -4494

$ws.Range("H126").Value = 2158.074
$ws.Range("I126").Value = 1875.2222
$ws.Range("J126").Value = 2723.7778
$ws.Range("K126").Value = 5625.6666
$ws.Range("L126").Value = 8171.3334
$ws.Range("M126").Value = -3155.6666
$ws.Range("N126").Value = -13111.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 610.5
$ws.Range("I97").Value = 622.5714
$ws.Range("J97").Value = 526
$ws.Range("K97").Value = 622.5714
$ws.Range("L97").Value = 526
$ws.Range("M97").Value = -126.5714
$ws.Range("N97").Value = -1518

$ws.Range("H113").Value = 13806.375
$ws.Range("I113").Value = 21050.2
$ws.Range("J113").Value = 1733.3334
$ws.Range("K113").Value = 21050.2
$ws.Range("L113").Value = 1733.3334
$ws.Range("M113").Value = -18880.2
$ws.Range("N113").Value = -6073.3334

$ws.Range("H122").Value = 2229.5715
$ws.Range("I122").Value = 1803.5
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 5410.5
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -2960.5
$ws.Range("N122").Value = -12100

$ws.Range("H126").Value = 1556.3334
$ws.Range("I126").Value = 1402.2
$ws.Range("J126").Value = 1749
$ws.Range("K126").Value = 4206.6
$ws.Range("L126").Value = 5247
$ws.Range("M126").Value = -1736.6
$ws.Range("N126").Value = -10187

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3081
$ws.Range("I40").Value = 3202.5
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3202.5
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3066.5
$ws.Range("N40").Value = -3272

$ws.Range("H46").Value = 875
$ws.Range("I46").Value = 566.6667
$ws.Range("J46").Value = 1800
$ws.Range("K46").Value = 566.6667
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -378.6667
$ws.Range("N46").Value = -2176

$ws.Range("H122").Value = 2599.8157
$ws.Range("I122").Value = 2393.5312
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 7180.5936
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -4730.5936
$ws.Range("N122").Value = -16000

$ws.Range("H132").Value = 1599.3016
$ws.Range("I132").Value = 854.1799999999999
$ws.Range("K132").Value = 2562.54
$ws.Range("M132").Value = -32.53999999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1381.75
$ws.Range("I132").Value = 1259.5652
$ws.Range("J132").Value = 1597.9231
$ws.Range("K132").Value = 3778.6956
$ws.Range("L132").Value = 4793.7693
$ws.Range("M132").Value = -1248.6956
$ws.Range("N132").Value = -9853.7693

$ws.Range("H136").Value = 1068.9535
$ws.Range("I136").Value = 860.65515
$ws.Range("J136").Value = 1500.4286
$ws.Range("K136").Value = 2581.96545
$ws.Range("L136").Value = 4501.2858
$ws.Range("M136").Value = -31.96545000000015
$ws.Range("N136").Value = -9601.2858
